$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.340.17'
$ws.Range("D3").Value = '3.284.41'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.58'
$ws.Range("D5").Style = $ws.Range("E5").Style
$ws.Range("E5").Value = '  +3.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.04'
$ws.Range("D6").Style = $ws.Range("E6").Style
$ws.Range("E6").Value = '  -3.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = $ws.Range("E7").Style
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.278.18'
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("E9").Value = '  -3.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("D10").Style = $ws.Range("E10").Style
$ws.Range("E10").Value = '  -5.93%  '
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.24'
$ws.Range("D12").Style = $ws.Range("E12").Style
$ws.Range("E12").Value = '  -2.77%  '
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("D14").Value = '3.811.03'
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '622.16'
$ws.Range("D15").Style = $ws.Range("E15").Style
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.37'
$ws.Range("D16").Style = $ws.Range("E16").Style
$ws.Range("E16").Value = '  -3.30%  '
$ws.Range("D17").Value = '65.484.89'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '3.286.22'
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("E20").Value = '  -2.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.89'
$ws.Range("D21").Style = $ws.Range("E21").Style
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("E22").Value = '  -2.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.91'
$ws.Range("D23").Style = $ws.Range("E23").Style
$ws.Range("E23").Value = '  -2.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.46'
$ws.Range("D24").Style = $ws.Range("E24").Style
$ws.Range("E24").Value = '  -2.76%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.33'
$ws.Range("D28").Style = $ws.Range("E28").Style
$ws.Range("E28").Value = '  -3.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.52'
$ws.Range("D29").Style = $ws.Range("E29").Style
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.32'
$ws.Range("D30").Style = $ws.Range("E30").Style
$ws.Range("E30").Value = '  -3.89%  '
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '567.99'
$ws.Range("D32").Style = $ws.Range("E32").Style
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.64'
$ws.Range("D33").Style = $ws.Range("E33").Style
$ws.Range("E33").Value = '  -10.47%  '
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("D35").Value = '3.835.67'
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("E36").Value = '  -2.04%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.40'
$ws.Range("D38").Style = $ws.Range("E38").Style
$ws.Range("E38").Value = '  -3.44%  '
$ws.Range("E39").Value = '  -2.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.42'
$ws.Range("D40").Style = $ws.Range("E40").Style
$ws.Range("E40").Value = '  +6.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.35'
$ws.Range("D41").Style = $ws.Range("E41").Style
$ws.Range("E41").Value = '  -4.88%  '
$ws.Range("D42").Value = '0.0₃0676'
$ws.Range("E42").Value = '  -6.96%  '
$ws.Range("E43").Value = '  -6.50%  '
$ws.Range("E44").Value = '  -5.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.329'
$ws.Range("D45").Style = $ws.Range("E45").Style
$ws.Range("E45").Value = '  -2.45%  '
$ws.Range("E46").Value = '  -4.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.00'
$ws.Range("D47").Style = $ws.Range("E47").Style
$ws.Range("E47").Value = '  -5.68%  '
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("E50").Value = '  -3.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.08'
$ws.Range("D51").Style = $ws.Range("E51").Style
$ws.Range("E51").Value = '  +4.97%  '
